$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Nmt_8"
$ws.Range("B9").Value = "p_186k+case"
$ws.Range("C9").Value = "m_20_100_1.1m"

$ws.Range("B9").Select()
